$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I0, IF), matching the look of
# the existing header row (bold font + thin border, centered/top-aligned)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Populate I2:J67 with the data values from the source diff
$data = @(
    @(7, 8),
    @(7, 8),
    @(7, 7),
    @(6, 6),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(6, 7),
    @(8, 8),
    @(8, 9),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(6, 7),
    @(11, 11),
    @(7, 7),
    @(8, 8),
    @(10, 10),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(7, 7),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(11, 12),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(10, 10),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(5, 5),
    @(7, 7),
    @(8, 8),
    @(4, 4),
    @(3, 3)
)

for ($k = 0; $k -lt $data.Count; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $data[$k][0]
    $ws.Cells.Item($row, 10).Value = $data[$k][1]
}
